$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.011094610950067
$ws.Range("C2").Value = 0.851912200597269
$ws.Range("D2").Value = 0.07773382499803461
$ws.Range("E2").Value = 0.05541757197566
$ws.Range("G2").Value = 0.002498352824656007
$ws.Range("I2").Value = 4.037666906858703
$ws.Range("M2").Value = 0.4877617833367225

$ws.Range("B3").Value = 0.9445053853133345
$ws.Range("C3").Value = 0.770090198141645
$ws.Range("D3").Value = 0.07057300904632768
$ws.Range("E3").Value = 0.05063322030689221
$ws.Range("G3").Value = 0.002506692214119027
$ws.Range("I3").Value = 3.701989948544707
$ws.Range("M3").Value = 0.4485988992850238

$ws.Range("B4").Value = 0.905138486006507
$ws.Range("C4").Value = 0.7203888840911645
$ws.Range("D4").Value = 0.06622831972570964
$ws.Range("E4").Value = 0.0477165122276908
$ws.Range("G4").Value = 0.002512059872913826
$ws.Range("I4").Value = 3.496110696767971
$ws.Range("M4").Value = 0.4249970262051761

$ws.Range("B5").Value = 0.889471005634789
$ws.Range("C5").Value = 0.7002655796318891
$ws.Range("D5").Value = 0.06447050239999896
$ws.Range("E5").Value = 0.04653293361720046
$ws.Range("G5").Value = 0.002514309712713722
$ws.Range("I5").Value = 3.412251916082198
$ws.Range("M5").Value = 0.4154880792310109

$ws.Range("B6").Value = 0.886891866654878
$ws.Range("C6").Value = 0.6969318559234239
$ws.Range("D6").Value = 0.0641793719211563
$ws.Range("E6").Value = 0.04633669619306247
$ws.Range("G6").Value = 0.002514687079209966
$ws.Range("I6").Value = 3.398329028553945
$ws.Range("M6").Value = 0.4139156314725341

$ws.Range("B7").Value = 0.9049256804531751
$ws.Range("C7").Value = 0.720116972558003
$ws.Range("D7").Value = 0.06620456239494388
$ws.Range("E7").Value = 0.04770053014411246
$ws.Range("G7").Value = 0.002512089961813152
$ws.Range("I7").Value = 3.494979615644695
$ws.Range("M7").Value = 0.4248683473931791

$ws.Range("B8").Value = 0.987814965539485
$ws.Range("C8").Value = 0.8235854063493662
$ws.Range("D8").Value = 0.07525371035582396
$ws.Range("E8").Value = 0.05376343436105202
$ws.Range("G8").Value = 0.002501177115772661
$ws.Range("I8").Value = 3.921866854210919
$ws.Range("M8").Value = 0.4741645433138189

$ws.Range("B9").Value = 1.162754261462908
$ws.Range("C9").Value = 1.030985521590253
$ws.Range("D9").Value = 0.0934323850919867
$ws.Range("E9").Value = 0.06583075442872399
$ws.Range("G9").Value = 0.002481724780382634
$ws.Range("I9").Value = 4.761684978658934
$ws.Range("M9").Value = 0.5744896393367469

$ws.Range("B10").Value = 1.29933389737397
$ws.Range("C10").Value = 1.186447642184078
$ws.Range("D10").Value = 0.1070821039570689
$ws.Range("E10").Value = 0.07482298068691762
$ws.Range("G10").Value = 0.002468600739639869
$ws.Range("I10").Value = 5.381639588208827
$ws.Range("M10").Value = 0.6506184819936038

$ws.Range("B11").Value = 1.363328137324402
$ws.Range("C11").Value = 1.257918685846732
$ws.Range("D11").Value = 0.1133622256748339
$ws.Range("E11").Value = 0.07894522485135269
$ws.Range("G11").Value = 0.002462879506383099
$ws.Range("I11").Value = 5.664602265979056
$ws.Range("M11").Value = 0.6858202454155418

$ws.Range("B12").Value = 1.387838292513152
$ws.Range("C12").Value = 1.285097052942717
$ws.Range("D12").Value = 0.1157510628113698
$ws.Range("E12").Value = 0.08051107880248054
$ws.Range("G12").Value = 0.002460748486064696
$ws.Range("I12").Value = 5.771911584220902
$ws.Range("M12").Value = 0.6992357593881309

$ws.Range("B13").Value = 1.3825471299985
$ws.Range("C13").Value = 1.279238546626459
$ws.Range("D13").Value = 0.1152361001992119
$ws.Range("E13").Value = 0.08017362359022684
$ws.Range("G13").Value = 0.002461205866027661
$ws.Range("I13").Value = 5.748793241771637
$ws.Range("M13").Value = 0.6963426371117549

$ws.Range("B14").Value = 1.36533900021152
$ws.Range("C14").Value = 1.260152350602255
$ws.Range("D14").Value = 0.1135585394921179
$ws.Range("E14").Value = 0.07907394965765491
$ws.Range("G14").Value = 0.002462703476827866
$ws.Range("I14").Value = 5.67342738376982
$ws.Range("M14").Value = 0.6869222173166634

$ws.Range("B15").Value = 1.354834858469644
$ws.Range("C15").Value = 1.248476507265309
$ws.Range("D15").Value = 0.1125323928729642
$ws.Range("E15").Value = 0.07840100822324558
$ws.Range("G15").Value = 0.002463625418634334
$ws.Range("I15").Value = 5.627284792310888
$ws.Range("M15").Value = 0.6811631571086281

$ws.Range("B16").Value = 1.295189980669022
$ws.Range("C16").Value = 1.181792445260101
$ws.Range("D16").Value = 0.106673151168863
$ws.Range("E16").Value = 0.07455424443099901
$ws.Range("G16").Value = 0.002468979618485356
$ws.Range("I16").Value = 5.363168140094899
$ws.Range("M16").Value = 0.6483297254813607

$ws.Range("B17").Value = 1.259083347304852
$ws.Range("C17").Value = 1.141080314101089
$ws.Range("D17").Value = 0.1030971875829323
$ws.Range("E17").Value = 0.0722026934179425
$ws.Range("G17").Value = 0.002472327788585017
$ws.Range("I17").Value = 5.201397937961559
$ws.Range("M17").Value = 0.6283359172022642

$ws.Range("B18").Value = 1.23849061335784
$ws.Range("C18").Value = 1.117734054281129
$ws.Range("D18").Value = 0.1010470201150184
$ws.Range("E18").Value = 0.0708531027210384
$ws.Range("G18").Value = 0.002474277019754757
$ws.Range("I18").Value = 5.108439333634323
$ws.Range("M18").Value = 0.616889486995106

$ws.Range("B19").Value = 1.231548029137514
$ws.Range("C19").Value = 1.109841305594159
$ws.Range("D19").Value = 0.100353992404834
$ws.Range("E19").Value = 0.07039665372545301
$ws.Range("G19").Value = 0.00247494103328849
$ws.Range("I19").Value = 5.076979486194034
$ws.Range("M19").Value = 0.6130230070944691

$ws.Range("B20").Value = 1.262908800114019
$ws.Range("C20").Value = 1.14540687388569
$ws.Range("D20").Value = 0.1034771649500783
$ws.Range("E20").Value = 0.07245271161409761
$ws.Range("G20").Value = 0.002471968945379749
$ws.Range("I20").Value = 5.218609455400696
$ws.Range("M20").Value = 0.6304587333450229

$ws.Range("B21").Value = 1.370385856564212
$ws.Range("C21").Value = 1.265755291679
$ws.Range("D21").Value = 0.1140509854298557
$ws.Range("E21").Value = 0.07939681650129415
$ws.Range("G21").Value = 0.002462262631669178
$ws.Range("I21").Value = 5.695559710759596
$ws.Range("M21").Value = 0.6896868788291499

$ws.Range("B22").Value = 1.442246822498532
$ws.Range("C22").Value = 1.345076318099416
$ws.Range("D22").Value = 0.121024153715112
$ws.Range("E22").Value = 0.08396360225697208
$ws.Range("G22").Value = 0.00245612568849126
$ws.Range("I22").Value = 6.008202873363928
$ws.Range("M22").Value = 0.728895127174269

$ws.Range("B23").Value = 1.403742144965349
$ws.Range("C23").Value = 1.302678222676093
$ws.Range("D23").Value = 0.1172965450499248
$ws.Range("E23").Value = 0.08152352182614209
$ws.Range("G23").Value = 0.002459382282386952
$ws.Range("I23").Value = 5.841246952518134
$ws.Range("M23").Value = 0.7079221333570871

$ws.Range("B24").Value = 1.261178798231583
$ws.Range("C24").Value = 1.143450650740533
$ws.Range("D24").Value = 0.1033053594608475
$ws.Range("E24").Value = 0.07233967109897321
$ws.Range("G24").Value = 0.002472131102599406
$ws.Range("I24").Value = 5.210827989112062
$ws.Range("M24").Value = 0.6294988580485921

$ws.Range("B25").Value = 1.114051328887626
$ws.Range("C25").Value = 0.974361795246125
$ws.Range("D25").Value = 0.08846516147946204
$ws.Range("E25").Value = 0.06254534745742291
$ws.Range("G25").Value = 0.002486780673932875
$ws.Range("I25").Value = 4.534071086306739
$ws.Range("M25").Value = 0.5469382274811494
